# Algos.xlsx edit:
#  - Add two new benchmark rows (algo 18 and 19, extractor "v3_small") to
#    the "MSH" sheet.
#  - Insert a new worksheet "error analysis(temp)" between "MSH" and "wiki",
#    populate it with the acronym-extractor-v3 failure analysis table, and
#    make it the active sheet.

$wb = $excel.ActiveWorkbook

$msh  = $wb.Worksheets.Item("MSH")
$wiki = $wb.Worksheets.Item("wiki")

# --- 1. MSH: add new benchmark rows for the v3_small extractor ------------

$msh.Cells.Item(22, 1).Value = 18
$msh.Cells.Item(22, 2).Value = 0.9092
$msh.Cells.Item(22, 3).Value = 0.007662
$msh.Cells.Item(22, 4).Value = "v3_small"
$msh.Cells.Item(22, 5).Value = "SVC"
$msh.Cells.Item(22, 6).Value = "all acronyms stored in uppercase"

$msh.Cells.Item(23, 1).Value = 19
$msh.Cells.Item(23, 4).Value = "v3_small"
$msh.Cells.Item(23, 5).Value = "LDA_multiclass"
$msh.Cells.Item(23, 6).Value = "all acronyms stored in uppercase"

# --- 2. Insert the new "error analysis(temp)" sheet right after MSH -------

$errSheet = $wb.Worksheets.Add($null, $msh)
$errSheet.Name = "error analysis(temp)"

# Column A: all 32 acronyms analysed in the v3 error-analysis pass.
$acronyms = @("TMP","CCL4","DI","HR","PCP","PAC","CAM","EPI","CCD","TAT","TEM","NBS","FAS","DAT","DDS","MBP","BSA","EM","TNT","TPA","TPO","ICE","BLM","TNC","CP","ADH","CDA","MCC","ALA","CDR","MAF","ORF")
$r = 2
foreach ($a in $acronyms) {
    $errSheet.Cells.Item($r, 1).Value = $a
    $r++
}

# Header row.
$errSheet.Cells.Item(1, 1).Value = "Acronym"
$errSheet.Cells.Item(1, 2).Value = "Errors"
$errSheet.Cells.Item(1, 3).Value = "Reason"

# Seed the four distinct failure-reason categories (first occurrences).
$errSheet.Cells.Item(3, 3).Value = "not detected as acronym"
$errSheet.Cells.Item(2, 3).Value = "case 1Capital, none"
$errSheet.Cells.Item(4, 3).Value = "case none capital"
$errSheet.Cells.Item(6, 3).Value = "1 capital case"

# Remaining columns B (error counts) and C (reason) for every row.
$details = @(
    @(2, 2, "case 1Capital, none"),
    @(3, 198, "not detected as acronym"),
    @(4, 18, "case none capital"),
    @(5, 9, "case none capital"),
    @(6, 1, "1 capital case"),
    @(7, 1, "case none capital"),
    @(8, 1, "1 capital case"),
    @(9, 59, "case 1Capital, none"),
    @(10, 1, "case 1Capital, none"),
    @(11, 159, "case 1Capital, none"),
    @(12, 2, "case none capital"),
    @(13, 1, "1 capital case"),
    @(14, 2, "case 1Capital, none"),
    @(15, 1, "case none capital"),
    @(16, 1, "1 capital case"),
    @(17, 6, "case 1Capital, none"),
    @(18, 1, "case none capital"),
    @(19, 1, "case none capital"),
    @(20, 1, "1 capital case"),
    @(21, 1, "1 capital case"),
    @(22, 11, "1 capital case"),
    @(23, 37, "case 1Capital, none"),
    @(24, 4, "case 1Capital, none"),
    @(25, 9, "1 capital case"),
    @(26, 13, "1 capital case"),
    @(27, 9, "case 1Capital, none"),
    @(28, 1, "case none capital"),
    @(29, 2, "1 capital case"),
    @(30, 97, "case 1Capital, none"),
    @(31, 2, "1 capital case"),
    @(32, 26, "case 1Capital, none"),
    @(33, 93, "case 1Capital, none")
)

foreach ($row in $details) {
    $rowNum = $row[0]
    $errSheet.Cells.Item($rowNum, 2).Value = $row[1]
    $errSheet.Cells.Item($rowNum, 3).Value = $row[2]
}

$errSheet.Columns.Item(3).ColumnWidth = 21.6328125

$msh.Range("B23").Select() | Out-Null

$errSheet.Activate() | Out-Null
$errSheet.Range("A2:C33").Select() | Out-Null
